$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1177.1852
$ws.Range("I28").Value = 979.2632
$ws.Range("K28").Value = 979.2632
$ws.Range("M28").Value = -494.2632
$ws.Range("H62").Value = 65917.3
$ws.Range("I62").Value = 205417.8
$ws.Range("J62").Value = 7792.0835
$ws.Range("K62").Value = 205417.8
$ws.Range("L62").Value = 7792.0835
$ws.Range("M62").Value = -204793.8
$ws.Range("N62").Value = -9040.083500000001
$ws.Range("H65").Value = 65917.3
$ws.Range("I65").Value = 205417.8
$ws.Range("J65").Value = 7792.0835
$ws.Range("K65").Value = 1027089
$ws.Range("L65").Value = 38960.4175
$ws.Range("M65").Value = -1023969
$ws.Range("N65").Value = -45200.4175
$ws.Range("H107").Value = 59722.883
$ws.Range("I107").Value = 84101.164
$ws.Range("K107").Value = 84101.164
$ws.Range("M107").Value = -82181.164
$ws.Range("H132").Value = 3157.1194
$ws.Range("I132").Value = 3156.377
$ws.Range("K132").Value = 9469.130999999999
$ws.Range("M132").Value = -6939.130999999999
$ws.Range("H137").Value = 254899.72
$ws.Range("I137").Value = 254899.72
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 764699.16
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -762149.16
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1558.931
$ws.Range("I2").Value = 1226.1904
$ws.Range("K2").Value = 1226.1904
$ws.Range("M2").Value = -1113.1904
$ws.Range("H32").Value = 1517.4639
$ws.Range("I32").Value = 1453.7094
$ws.Range("K32").Value = 1453.7094
$ws.Range("M32").Value = -1166.7094
$ws.Range("H61").Value = 3842.5
$ws.Range("I61").Value = 3842.5
$ws.Range("K61").Value = 3842.5
$ws.Range("M61").Value = -3630.5
$ws.Range("H74").Value = 22926.234
$ws.Range("I74").Value = 4505.0303
$ws.Range("K74").Value = 4505.0303
$ws.Range("M74").Value = -3631.0303
$ws.Range("H77").Value = 22926.234
$ws.Range("I77").Value = 4505.0303
$ws.Range("K77").Value = 22525.1515
$ws.Range("M77").Value = -18157.1515
$ws.Range("H97").Value = 10066.913
$ws.Range("I97").Value = 6990.294
$ws.Range("K97").Value = 6990.294
$ws.Range("M97").Value = -6494.294
$ws.Range("H116").Value = 1558.931
$ws.Range("I116").Value = 1226.1904
$ws.Range("K116").Value = 1226.1904
$ws.Range("M116").Value = 1067.8096
$ws.Range("H132").Value = 6197.5557
$ws.Range("I132").Value = 3995
$ws.Range("J132").Value = 6472.875
$ws.Range("K132").Value = 11985
$ws.Range("L132").Value = 19418.625
$ws.Range("M132").Value = -9455
$ws.Range("N132").Value = -24478.625
$ws.Range("H136").Value = 3842.5
$ws.Range("I136").Value = 3842.5
$ws.Range("K136").Value = 11527.5
$ws.Range("M136").Value = -8977.5
$ws.Range("H139").Value = 241722
$ws.Range("J139").Value = 241722
$ws.Range("L139").Value = 241722
$ws.Range("N139").Value = -252002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1558.931
$ws.Range("I3").Value = 1226.1904
$ws.Range("K3").Value = 1226.1904
$ws.Range("M3").Value = -1112.1904
$ws.Range("H134").Value = 3695.0293
$ws.Range("I134").Value = 1688.3077
$ws.Range("J134").Value = 10216.875
$ws.Range("K134").Value = 5064.9231
$ws.Range("L134").Value = 30650.625
$ws.Range("M134").Value = -2529.9231
$ws.Range("N134").Value = -35720.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1618.88
$ws.Range("I16").Value = 1493.1305
$ws.Range("K16").Value = 1493.1305
$ws.Range("M16").Value = -1206.1305
$ws.Range("H31").Value = 42361.81
$ws.Range("I31").Value = 2245.7693
$ws.Range("K31").Value = 2245.7693
$ws.Range("M31").Value = -1950.7693
$ws.Range("H34").Value = 42361.81
$ws.Range("I34").Value = 2245.7693
$ws.Range("K34").Value = 2245.7693
$ws.Range("M34").Value = -2043.7693
$ws.Range("H58").Value = 2421.0454
$ws.Range("I58").Value = 2351.1875
$ws.Range("J58").Value = 2607.3333
$ws.Range("K58").Value = 2351.1875
$ws.Range("L58").Value = 2607.3333
$ws.Range("M58").Value = -2148.1875
$ws.Range("N58").Value = -3013.3333
$ws.Range("H94").Value = 4666.75
$ws.Range("J94").Value = 5631.5
$ws.Range("L94").Value = 5631.5
$ws.Range("N94").Value = -6533.5
$ws.Range("H113").Value = 1618.88
$ws.Range("I113").Value = 1493.1305
$ws.Range("K113").Value = 1493.1305
$ws.Range("M113").Value = 676.8695
$ws.Range("H132").Value = 25779.34
$ws.Range("I132").Value = 36837.668
$ws.Range("K132").Value = 110513.004
$ws.Range("M132").Value = -107983.004
$ws.Range("H134").Value = 2950.2334
$ws.Range("I134").Value = 2139.8948
$ws.Range("K134").Value = 6419.6844
$ws.Range("M134").Value = -3884.6844
$ws.Range("H136").Value = 2421.0454
$ws.Range("I136").Value = 2351.1875
$ws.Range("J136").Value = 2607.3333
$ws.Range("K136").Value = 7053.5625
$ws.Range("L136").Value = 7821.999899999999
$ws.Range("M136").Value = -4503.5625
$ws.Range("N136").Value = -12921.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 9685461
$ws.Range("I4").Value = 11033280
$ws.Range("K4").Value = 33099840
$ws.Range("M4").Value = -33099728
$ws.Range("H34").Value = 596.6875
$ws.Range("I34").Value = 179.45454
$ws.Range("J34").Value = 1514.6
$ws.Range("K34").Value = 538.3636200000001
$ws.Range("L34").Value = 4543.799999999999
$ws.Range("M34").Value = -454.3636200000001
$ws.Range("N34").Value = -4711.799999999999
$ws.Range("H39").Value = 6374.5
$ws.Range("J39").Value = 10000
$ws.Range("L39").Value = 30000
$ws.Range("N39").Value = -30588
$ws.Range("H55").Value = 65287.5
$ws.Range("J55").Value = 129875
$ws.Range("L55").Value = 389625
$ws.Range("N55").Value = -389979
$ws.Range("H76").Value = 1000
$ws.Range("I76").Value = 1000
$ws.Range("K76").Value = 3000
$ws.Range("M76").Value = -2617
$ws.Range("H79").Value = 1000
$ws.Range("I79").Value = 1000
$ws.Range("K79").Value = 3000
$ws.Range("M79").Value = -1674
$ws.Range("H132").Value = 1318.36
$ws.Range("I132").Value = 1286.2354
$ws.Range("J132").Value = 1386.625
$ws.Range("K132").Value = 11576.1186
$ws.Range("L132").Value = 12479.625
$ws.Range("M132").Value = -9046.1186
$ws.Range("N132").Value = -17539.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 77117.60000000001
$ws.Range("J45").Value = 77117.60000000001
$ws.Range("L45").Value = 77117.60000000001
$ws.Range("N45").Value = -78235.60000000001
$ws.Range("H47").Value = 20666.666
$ws.Range("J47").Value = 20666.666
$ws.Range("L47").Value = 20666.666
$ws.Range("N47").Value = -21802.666
$ws.Range("H51").Value = 81727.27
$ws.Range("J51").Value = 81727.27
$ws.Range("L51").Value = 81727.27
$ws.Range("N51").Value = -82745.27
$ws.Range("H122").Value = 564610.75
$ws.Range("I122").Value = 749506.7
$ws.Range("K122").Value = 2248520.1
$ws.Range("M122").Value = -2246070.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2158.9355
$ws.Range("I132").Value = 1851.4546
$ws.Range("K132").Value = 5554.3638
$ws.Range("M132").Value = -3024.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 130032.55
$ws.Range("I132").Value = 1882.1034
$ws.Range("J132").Value = 395487.03
$ws.Range("K132").Value = 5646.3102
$ws.Range("L132").Value = 1186461.09
$ws.Range("M132").Value = -3116.3102
$ws.Range("N132").Value = -1191521.09
$ws.Range("H136").Value = 9510.576999999999
$ws.Range("I136").Value = 10585.228
$ws.Range("J136").Value = 3600
$ws.Range("K136").Value = 31755.684
$ws.Range("L136").Value = 10800
$ws.Range("M136").Value = -29205.684
$ws.Range("N136").Value = -15900
